$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 161, pushing the existing rows 161:245 down to 162:246
$ws.Rows("161:161").Insert()

# Populate the newly inserted row 161 with the new record.
# Columns A,B,C,E,F,G,H,I,R mirror the surrounding "Albahaca" dataset rows,
# only D,J,K,L,M,N,O,P,Q carry the new values from the diff.
$ws.Cells.Item(161, 1).Value = 9
$ws.Cells.Item(161, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(161, 3).Value = "Metropolitana"
$ws.Cells.Item(161, 4).Value = 44529
$ws.Cells.Item(161, 5).Value = 13
$ws.Cells.Item(161, 6).Value = 100112052
$ws.Cells.Item(161, 7).Value = "Albahaca"
$ws.Cells.Item(161, 8).Value = "Sin especificar"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 106
$ws.Cells.Item(161, 11).Value = 3500
$ws.Cells.Item(161, 12).Value = 4000
$ws.Cells.Item(161, 13).Value = 3750
$ws.Cells.Item(161, 14).Value = "$/paquete"
$ws.Cells.Item(161, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(161, 16).Value = 3750
$ws.Cells.Item(161, 17).Value = 1
$ws.Cells.Item(161, 18).Value = "Hortaliza"
